# Swap the contents of columns C and D (codeforiati:group-code and
# codeforiati:group-name, including the header row) for every used row
# on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
